$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / cell A2 gets the new text (shared string)
$ws.Range("A2").Value = "Deuxième ligne modifié"

# Column A width matches the target width
$ws.Columns("A").ColumnWidth = 26.7109375

# Leave the cursor/selection on B2, as in the saved file
$ws.Range("B2").Select()
